$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts the string into a
# numeric value (losing the exact original text representation).
$textCells = @("D5", "D6", "D8", "D9", "D11", "D12", "D14", "D16", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.947.43'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.643.26'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '215.51'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").Value = '0.5087'
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '0.2567'
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").Value = '0.06395'
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("D11").Value = '0.07775'
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").Value = '4.303'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '1.642.27'
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").Value = '0.5457'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").Value = '0.0₅7858'
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").Value = '64.70'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").Value = '25.996.25'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '197.51'
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").Value = '4.433'
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").Value = '9.967'
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '6.039'
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("D24").Value = '1.873'
$ws.Range("E24").Value = '  -3.02%  '
$ws.Range("D25").Value = '141.64'
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").Value = '0.1144'
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("D27").Value = '6.893'
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("D28").Value = '15.73'
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").Value = '1.238'
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("D30").Value = '0.05038'
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("D31").Value = '3.263'
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("D32").Value = '3.190'
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").Value = '1.544'
$ws.Range("E33").Value = '  -0.35%  '
$ws.Range("D34").Value = '2.364'
$ws.Range("E34").Value = '  -1.07%  '
$ws.Range("D35").Value = '0.8952'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").Value = '2.586'
$ws.Range("E36").Value = '  -1.93%  '
$ws.Range("D37").Value = '1.130.29'
$ws.Range("E37").Value = '  -3.76%  '
$ws.Range("D38").Value = '0.5515'
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("D39").Value = '0.01554'
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").Value = '1.006'
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("E41").Value = '  +18.86%  '
$ws.Range("D42").Value = '2.546'
$ws.Range("E42").Value = '  -1.27%  '
$ws.Range("D43").Value = '5.642'
$ws.Range("E43").Value = '  -1.23%  '
$ws.Range("D44").Value = '0.8176'
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("D45").Value = '100.09'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '1.778.46'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").Value = '54.92'
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("E51").Value = '  -0.33%  '
